$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B11 currently holds the literal text "R40" (rule name for the last row of
# the rules table). The target edit renames that rule label to the literal
# text "1" -- still a text value (not a number) stored in the same cell
# style/format as before.
#
# A plain `$ws.Range("B11").Value = "1"` would have Excel auto-detect the
# numeric-looking string and store it as a Number, changing the cell's type.
# To keep it as text (matching the original t="s" shared-string cell type)
# without disturbing B11's existing style, format a scratch cell as Text,
# write "1" into it there, then copy/paste-values that text into B11 - this
# mirrors how a pre-formatted-as-text source cell copies as text without
# touching the destination's own number format.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
$scratch.Clear()
